$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.354516744613647
$ws.Range("B1").Value = 2.643781185150146
$ws.Range("C1").Value = 2.003476619720459
$ws.Range("D1").Value = 1.856452465057373
$ws.Range("E1").Value = 1.893656611442566
